$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 156140
$ws.Range("C4").Value = 147232
$ws.Range("C5").Value = 8908
$ws.Range("C8").Value = 63.65
